# Add 2022-Q3 data:
#  - record the new quarter's totals on the "总计" sheet (pushing the old
#    2020-Q4 total row down to row 3)
#  - insert a brand-new "2022-Q3" worksheet (positioned between "总计" and
#    "2020-Q4") holding the per-fund position breakdown for the quarter.

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item("总计")
$q4sheet = $wb.Worksheets.Item("2020-Q4")

function Set-TextValue($range, [string]$text) {
    # Force a string like "011189" or "17.72" to be written verbatim as
    # text (no leading-zero / trailing-zero loss from numeric coercion),
    # then drop the temporary "@" text format so the cell is left with no
    # explicit style (matching freshly authored, unstyled data cells).
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# --- 1. Update the "总计" (summary) sheet -------------------------------
# Row 2 used to hold the 2020-Q4 totals; it now becomes the 2022-Q3 totals,
# and the old 2020-Q4 totals move down to row 3 (copy A2's formatting first
# so the new "index" cell keeps the same bordered/bold style as A2).
$summary.Range("A2").Copy($summary.Range("A3"))
$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2020-Q4"
$summary.Range("C3").Value = 2
$summary.Range("D3").Value = 0.02

$summary.Range("B2").Value = "2022-Q3"
$summary.Range("D2").Value = 0.46

# --- 2. Insert the new "2022-Q3" sheet ----------------------------------
# Placed right after "总计" and before "2020-Q4", matching the new sheet
# order: 总计, 2022-Q3, 2020-Q4.
$newSheet = $wb.Worksheets.Add($q4sheet)
$newSheet.Name = "2022-Q3"

# Header row, copied style-and-all from the "总计" header so it shares the
# same bold/bordered look (style index) as the rest of the workbook.
$summary.Range("B1").Copy($newSheet.Range("B1"))
$newSheet.Range("B1").Value = "基金代码"
$summary.Range("C1").Copy($newSheet.Range("C1"))
$newSheet.Range("C1").Value = "基金名称"
$summary.Range("D1").Copy($newSheet.Range("D1"))
$newSheet.Range("D1").Value = "基金规模"
$summary.Range("B1").Copy($newSheet.Range("E1"))
$newSheet.Range("E1").Value = "股票总仓位"
$summary.Range("B1").Copy($newSheet.Range("F1"))
$newSheet.Range("F1").Value = "仓位占比"
$summary.Range("D1").Copy($newSheet.Range("G1"))
$newSheet.Range("G1").Value = "持有市值(亿元)"
$summary.Range("B1").Copy($newSheet.Range("H1"))
$newSheet.Range("H1").Value = "仓位排名"

# Row 2 (index cell A2 copies the "总计" A2 style too)
$summary.Range("A2").Copy($newSheet.Range("A2"))
$newSheet.Range("A2").Value = 0
Set-TextValue $newSheet.Range("B2") "011189"
$newSheet.Range("C2").Value = "建信智汇优选一年持有期混合（MOM）"
Set-TextValue $newSheet.Range("D2") "17.72"
Set-TextValue $newSheet.Range("E2") "54.61"
Set-TextValue $newSheet.Range("F2") "1.58"
Set-TextValue $newSheet.Range("G2") "0.2800"
$newSheet.Range("H2").Value = 8

# Row 3
$summary.Range("A2").Copy($newSheet.Range("A3"))
$newSheet.Range("A3").Value = 1
Set-TextValue $newSheet.Range("B3") "011481"
$newSheet.Range("C3").Value = "广发瑞锦一年定开混合"
Set-TextValue $newSheet.Range("D3") "2.62"
Set-TextValue $newSheet.Range("E3") "89.31"
Set-TextValue $newSheet.Range("F3") "7.00"
Set-TextValue $newSheet.Range("G3") "0.1834"
$newSheet.Range("H3").Value = 2

Write-Host "Sheets now:"
foreach ($ws in $wb.Worksheets) {
    Write-Host (" - {0}" -f $ws.Name)
}
